$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 252
$ws.Range("I2").Value = 601
$ws.Range("J2").Value = 2695
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 730
$ws.Range("M2").Value = 46
$ws.Range("N2").Value = 447
$ws.Range("O2").Value = 1
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 285
$ws.Range("T2").Value = 449
$ws.Range("U2").Value = 44
$ws.Range("V2").Value = 4105
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 4215
$ws.Range("Z2").Value = 76
$ws.Range("AA2").Value = 23
